# Add a new slide at the end of the deck: "StringUtils from Apache"
$p = $ppt.ActivePresentation

$newIndex = $p.Slides.Count + 1
$slide = $p.Slides.Add($newIndex, 2)  # ppLayoutText -> "Title and Content" (same layout as the rest of the deck)

# Title placeholder: "StringUtils" + " from Apache" (two runs)
$titleRange = $slide.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "StringUtils"
$inserted = $titleRange.InsertAfter(" from Apache")

# Body placeholder: bold code snippet
$bodyRange = $slide.Shapes.Item(2).TextFrame.TextRange
$bodyRange.Text = "import org.apache.commons.lang3.StringUtils;"
$bodyRange.Font.Bold = $true
